# Cheat Sheet.xlsx update:
#   - Insert a new row for the "ln -s" entry group so it also documents
#     `readlink -f link` (get a symbolic link destination).
#   - This pushes the "Volumes/File System" / "Variables" section (and
#     everything below it) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 13 (pushes old row 13.. down to 14..)
$ws.Rows(13).Insert()

# Match the row height used by the other data rows in the sheet.
$ws.Rows(13).RowHeight = 18.75

# Fill in the new entry in columns C/D (A/B stay blank on this row).
$ws.Range("C13").Value = "readlink -f link"
$ws.Range("D13").Value = "get a symbolic link destination"

# Cosmetic: keep the "unzip" row's height consistent with the target layout
# after the shift caused by the row insertion above.
$ws.Rows(25).RowHeight = 18.75
